# Scheduled pricing-refresh runner: update market price / leve profit
# figures (columns H-N) on several rows across the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 556.5
$ws.Range("I2").Value = 556.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 556.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -443.5
$ws.Range("N2").ClearContents()
$ws.Range("H17").Value = 2214.76
$ws.Range("J17").Value = 2214.76
$ws.Range("L17").Value = 6644.280000000001
$ws.Range("N17").Value = -6980.280000000001
$ws.Range("H51").Value = 13749.25
$ws.Range("I51").Value = 17998.5
$ws.Range("J51").Value = 9500
$ws.Range("K51").Value = 17998.5
$ws.Range("L51").Value = 9500
$ws.Range("M51").Value = -17514.5
$ws.Range("N51").Value = -10468
$ws.Range("H87").Value = 85712.86
$ws.Range("J87").Value = 85712.86
$ws.Range("L87").Value = 85712.86
$ws.Range("N87").Value = -88208.86
$ws.Range("H90").Value = 85712.86
$ws.Range("J90").Value = 85712.86
$ws.Range("L90").Value = 257138.58
$ws.Range("N90").Value = -269618.58
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 1042
$ws.Range("I141").Value = 1042
$ws.Range("K141").Value = 3126
$ws.Range("M141").Value = 2054

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 11003.5
$ws.Range("I25").Value = 2007
$ws.Range("K25").Value = 2007
$ws.Range("M25").Value = -1605
$ws.Range("H32").Value = 3129.1853
$ws.Range("I32").Value = 1983.56
$ws.Range("K32").Value = 1983.56
$ws.Range("M32").Value = -1696.56
$ws.Range("H110").Value = 5233.273
$ws.Range("I110").Value = 5656.6
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 5656.6
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = -3611.6
$ws.Range("N110").Value = -5090
$ws.Range("H131").Value = 69885
$ws.Range("J131").Value = 69885
$ws.Range("L131").Value = 69885
$ws.Range("N131").Value = -79965
$ws.Range("H135").Value = 86142.336
$ws.Range("J135").Value = 86142.336
$ws.Range("L135").Value = 86142.336
$ws.Range("N135").Value = -96282.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3483
$ws.Range("I20").Value = 3814.5557
$ws.Range("J20").Value = 3056.7144
$ws.Range("K20").Value = 3814.5557
$ws.Range("L20").Value = 3056.7144
$ws.Range("M20").Value = -3567.5557
$ws.Range("N20").Value = -3550.7144
$ws.Range("H86").Value = 2621.4443
$ws.Range("I86").Value = 2842
$ws.Range("K86").Value = 2842
$ws.Range("M86").Value = -1719
$ws.Range("H89").Value = 2621.4443
$ws.Range("I89").Value = 2842
$ws.Range("K89").Value = 14210
$ws.Range("M89").Value = -8594
$ws.Range("H99").Value = 2603

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2601.3333
$ws.Range("I62").Value = 2601.3333
$ws.Range("K62").Value = 2601.3333
$ws.Range("M62").Value = -1977.3333
$ws.Range("H65").Value = 2601.3333
$ws.Range("I65").Value = 2601.3333
$ws.Range("K65").Value = 13006.6665
$ws.Range("M65").Value = -9886.666499999999
$ws.Range("H107").Value = 1081.909
$ws.Range("I107").Value = 655.6667
$ws.Range("K107").Value = 655.6667
$ws.Range("M107").Value = 1264.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1448.3334
$ws.Range("I5").Value = 1422.5
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 4267.5
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = -4155.5
$ws.Range("N5").Value = -4724
$ws.Range("H33").Value = 1000
$ws.Range("I33").Value = 1000
$ws.Range("K33").Value = 6000
$ws.Range("M33").Value = -5717
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H110").Value = 8888
$ws.Range("I110").Value = 8888
$ws.Range("K110").Value = 26664
$ws.Range("M110").Value = -22574
$ws.Range("H121").Value = 1429069.6
$ws.Range("I121").Value = 581.1667
$ws.Range("J121").Value = 10000000
$ws.Range("K121").Value = 1743.5001
$ws.Range("L121").Value = 30000000
$ws.Range("M121").Value = -433.5001
$ws.Range("N121").Value = -30002620
$ws.Range("H122").Value = 586.25
$ws.Range("I122").Value = 586.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5276.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2826.25
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 15500
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H135").Value = 1448.3334
$ws.Range("I135").Value = 1422.5
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 12802.5
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -10267.5
$ws.Range("N135").Value = -18570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 7000
$ws.Range("J22").Value = 7000
$ws.Range("L22").Value = 7000
$ws.Range("N22").Value = -8058
$ws.Range("H122").Value = 1476.8125
$ws.Range("I122").Value = 1476.8125
$ws.Range("K122").Value = 4430.4375
$ws.Range("M122").Value = -1980.4375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("H82").Value = 2718.5386
$ws.Range("I82").Value = 3776
$ws.Range("K82").Value = 3776
$ws.Range("M82").Value = -3415
$ws.Range("H85").Value = 2718.5386
$ws.Range("I85").Value = 3776
$ws.Range("K85").Value = 3776
$ws.Range("M85").Value = -2528
$ws.Range("H136").Value = 21741002
$ws.Range("J136").Value = 100002240
$ws.Range("L136").Value = 300006720
$ws.Range("N136").Value = -300011820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 20000
$ws.Range("K20").Value = 20000
$ws.Range("M20").Value = -19760
$ws.Range("H45").Value = 129999
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 129999
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 129999
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -130981
$ws.Range("H132").Value = 4419.552
$ws.Range("J132").Value = 11500
$ws.Range("L132").Value = 34500
$ws.Range("N132").Value = -39560
